$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-18 08:35:11"
$wsZh.Range("H3").Value = "2016-03-18 08:35:29"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-18 08:35:14"
$wsDe.Range("H3").Value = "2016-03-18 08:35:34"
